# Apply the edit described by the diff:
#  - delete the old "2005 industry breakdown" data block (rows 121-131), which shifts
#    all subsequent rows (footnotes / source lines / hyperlink) up by 11
#  - fix the hyperlink anchor so it still points at the (now shifted) source-link cell,
#    restoring the original font/formatting that Hyperlinks.Add() overwrites by default
#  - update the footnote text so "Erhebungsjahres 2005" becomes "Erhebungsjahres 2006"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# 1) Delete rows 121-131 (entire old 2005 data table + its "Insgesamt" total row)
$ws.Range("A121:A131").EntireRow.Delete() | Out-Null

# 2) Re-anchor the hyperlink that used to live on B143 and now lives on B132
$ws.Range("B143").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B132"), "https://www.integrationsmonitoring.niedersachsen.de/") | Out-Null

# Hyperlinks.Add() resets the cell's font to the generic default "Hyperlink" style
# (Calibri 11); restore the sheet's original small-footnote hyperlink font so the
# cell keeps looking the way it did (and reuses the same cell style as before).
$linkCell = $ws.Range("B132")
$linkCell.Font.Name = "NDSFrutiger 45 Light"
$linkCell.Font.Size = 6

# 3) Update footnote 1) text: the referenced survey year changes from 2005 to 2006
$footnoteCell = $ws.Range("B122")
$oldText = $footnoteCell.Text
$newText = $oldText.Replace("Erhebungsjahres 2005", "Erhebungsjahres 2006")
$footnoteCell.Value = $newText
